$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.043840238251543
$ws.Cells.Item(2, 4).Value = 1.050778332049941
$ws.Cells.Item(2, 5).Value = 1.057197741469932
$ws.Cells.Item(2, 6).Value = 1.064100322796586
$ws.Cells.Item(2, 9).Value = 1.044201513436368
$ws.Cells.Item(2, 10).Value = 1.048908528474003
$ws.Cells.Item(2, 11).Value = 1.053531664543123
$ws.Cells.Item(2, 12).Value = 1.059933369664783
$ws.Cells.Item(2, 13).Value = 1.066817172114993
$ws.Cells.Item(2, 14).Value = 1.020183794367798

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.044640151352149
$ws.Cells.Item(3, 4).Value = 1.051392601602742
$ws.Cells.Item(3, 5).Value = 1.057948172177489
$ws.Cells.Item(3, 6).Value = 1.064852070172603
$ws.Cells.Item(3, 9).Value = 1.044379729686881
$ws.Cells.Item(3, 10).Value = 1.049356194338291
$ws.Cells.Item(3, 11).Value = 1.053958959262137
$ws.Cells.Item(3, 12).Value = 1.0604977503138
$ws.Cells.Item(3, 13).Value = 1.06738421785964
$ws.Cells.Item(3, 14).Value = 1.020333860296197

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.045158566516036
$ws.Cells.Item(4, 4).Value = 1.051790790611446
$ws.Cells.Item(4, 5).Value = 1.058434879781575
$ws.Cells.Item(4, 6).Value = 1.065339595425865
$ws.Cells.Item(4, 9).Value = 1.044494298542021
$ws.Cells.Item(4, 10).Value = 1.049645985069257
$ws.Cells.Item(4, 11).Value = 1.05423546714439
$ws.Cells.Item(4, 12).Value = 1.060863402253985
$ws.Cells.Item(4, 13).Value = 1.067751559487011
$ws.Cells.Item(4, 14).Value = 1.02043096728422

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.045376702142348
$ws.Cells.Item(5, 4).Value = 1.0519583586354
$ws.Cells.Item(5, 5).Value = 1.058639759775004
$ws.Cells.Item(5, 6).Value = 1.065544810900784
$ws.Cells.Item(5, 9).Value = 1.044542283199973
$ws.Cells.Item(5, 10).Value = 1.049767840671374
$ws.Cells.Item(5, 11).Value = 1.054351714349374
$ws.Cells.Item(5, 12).Value = 1.061017230747173
$ws.Cells.Item(5, 13).Value = 1.06790608994837
$ws.Cells.Item(5, 14).Value = 1.020471791583855

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.0454133393881
$ws.Cells.Item(6, 4).Value = 1.051986503919262
$ws.Cells.Item(6, 5).Value = 1.058674175688529
$ws.Cells.Item(6, 6).Value = 1.065579282657662
$ws.Cells.Item(6, 9).Value = 1.044550329451982
$ws.Cells.Item(6, 10).Value = 1.049788302340137
$ws.Cells.Item(6, 11).Value = 1.054371232918679
$ws.Cells.Item(6, 12).Value = 1.061043065536605
$ws.Cells.Item(6, 13).Value = 1.067932042108896
$ws.Cells.Item(6, 14).Value = 1.020478646179253

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.045161480495485
$ws.Cells.Item(7, 4).Value = 1.051793029000531
$ws.Cells.Item(7, 5).Value = 1.058437616347211
$ws.Cells.Item(7, 6).Value = 1.065342336506723
$ws.Cells.Item(7, 9).Value = 1.04449494042397
$ws.Cells.Item(7, 10).Value = 1.049647613203146
$ws.Cells.Item(7, 11).Value = 1.054237020434066
$ws.Cells.Item(7, 12).Value = 1.060865457293953
$ws.Cells.Item(7, 13).Value = 1.067753623939679
$ws.Cells.Item(7, 14).Value = 1.020431512779095

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.044110402125336
$ws.Cells.Item(8, 4).Value = 1.05098577800255
$ws.Cells.Item(8, 5).Value = 1.057451118056409
$ws.Cells.Item(8, 6).Value = 1.06435415143951
$ws.Cells.Item(8, 9).Value = 1.044261897062143
$ws.Cells.Item(8, 10).Value = 1.049059793366182
$ws.Cells.Item(8, 11).Value = 1.05367606588407
$ws.Cells.Item(8, 12).Value = 1.060124008502753
$ws.Cells.Item(8, 13).Value = 1.067008718815739
$ws.Cells.Item(8, 14).Value = 1.020234508656313

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.042264616724049
$ws.Cells.Item(9, 4).Value = 1.049568864305601
$ws.Cells.Item(9, 5).Value = 1.055721520291217
$ws.Cells.Item(9, 6).Value = 1.062621320300943
$ws.Cells.Item(9, 9).Value = 1.043845543310576
$ws.Cells.Item(9, 10).Value = 1.048024967561291
$ws.Cells.Item(9, 11).Value = 1.052687804888121
$ws.Cells.Item(9, 12).Value = 1.05882107815834
$ws.Cells.Item(9, 13).Value = 1.065699433297098
$ws.Cells.Item(9, 14).Value = 1.019887417558026

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041038476504109
$ws.Cells.Item(10, 4).Value = 1.048628113067029
$ws.Cells.Item(10, 5).Value = 1.054574456557743
$ws.Cells.Item(10, 6).Value = 1.061471927936815
$ws.Cells.Item(10, 9).Value = 1.04356418768241
$ws.Cells.Item(10, 10).Value = 1.047335830751951
$ws.Cells.Item(10, 11).Value = 1.052029191425445
$ws.Cells.Item(10, 12).Value = 1.057954977545814
$ws.Cells.Item(10, 13).Value = 1.064828920550193
$ws.Cells.Item(10, 14).Value = 1.019656091023506

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040508608213436
$ws.Cells.Item(11, 4).Value = 1.048221697023877
$ws.Cells.Item(11, 5).Value = 1.05407921557349
$ws.Cells.Item(11, 6).Value = 1.060975638112973
$ws.Cells.Item(11, 9).Value = 1.043441469440606
$ws.Cells.Item(11, 10).Value = 1.047037621667622
$ws.Cells.Item(11, 11).Value = 1.051744076283419
$ws.Cells.Item(11, 12).Value = 1.057580565480866
$ws.Cells.Item(11, 13).Value = 1.064452556946544
$ws.Cells.Item(11, 14).Value = 1.019555946525114

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.040311952537484
$ws.Cells.Item(12, 4).Value = 1.048070878557623
$ws.Cells.Item(12, 5).Value = 1.053895480333933
$ws.Cells.Item(12, 6).Value = 1.060791507241976
$ws.Cells.Item(12, 9).Value = 1.043395753595333
$ws.Cells.Item(12, 10).Value = 1.046926883761787
$ws.Cells.Item(12, 11).Value = 1.051638183615591
$ws.Cells.Item(12, 12).Value = 1.057441586342976
$ws.Cells.Item(12, 13).Value = 1.064312846797376
$ws.Cells.Item(12, 14).Value = 1.019518752138413

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.04035412853032
$ws.Cells.Item(13, 4).Value = 1.04810322314782
$ws.Cells.Item(13, 5).Value = 1.053934882207019
$ws.Cells.Item(13, 6).Value = 1.06083099425224
$ws.Cells.Item(13, 9).Value = 1.043405565797203
$ws.Cells.Item(13, 10).Value = 1.046950636027841
$ws.Cells.Item(13, 11).Value = 1.051660897395487
$ws.Cells.Item(13, 12).Value = 1.05747139355027
$ws.Cells.Item(13, 13).Value = 1.064342811085421
$ws.Cells.Item(13, 14).Value = 1.019526730283062

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040492349298172
$ws.Cells.Item(14, 4).Value = 1.048209227402881
$ws.Cells.Item(14, 5).Value = 1.054064023474701
$ws.Cells.Item(14, 6).Value = 1.060960413434249
$ws.Cells.Item(14, 9).Value = 1.043437693259581
$ws.Cells.Item(14, 10).Value = 1.047028467412943
$ws.Cells.Item(14, 11).Value = 1.051735322913557
$ws.Cells.Item(14, 12).Value = 1.057569075490772
$ws.Cells.Item(14, 13).Value = 1.064441006655475
$ws.Cells.Item(14, 14).Value = 1.019552871946221

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.04057753306807
$ws.Cells.Item(15, 4).Value = 1.04827455907195
$ws.Cells.Item(15, 5).Value = 1.054143620798446
$ws.Cells.Item(15, 6).Value = 1.06104018119035
$ws.Cells.Item(15, 9).Value = 1.043457470472559
$ws.Cells.Item(15, 10).Value = 1.047076425945896
$ws.Cells.Item(15, 11).Value = 1.051781180537308
$ws.Cells.Item(15, 12).Value = 1.057629273089322
$ws.Cells.Item(15, 13).Value = 1.064501519908148
$ws.Cells.Item(15, 14).Value = 1.019568979198702

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041073664587257
$ws.Cells.Item(16, 4).Value = 1.048655105441575
$ws.Cells.Item(16, 5).Value = 1.054607354706169
$ws.Cells.Item(16, 6).Value = 1.061504894843087
$ws.Cells.Item(16, 9).Value = 1.043572313417247
$ws.Cells.Item(16, 10).Value = 1.047355626067315
$ws.Cells.Item(16, 11).Value = 1.052048115152379
$ws.Cells.Item(16, 12).Value = 1.057979839152292
$ws.Cells.Item(16, 13).Value = 1.064853910812738
$ws.Cells.Item(16, 14).Value = 1.019662737779791

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.041385159565942
$ws.Cells.Item(17, 4).Value = 1.048894063989887
$ws.Cells.Item(17, 5).Value = 1.054898631197909
$ws.Cells.Item(17, 6).Value = 1.06179677508314
$ws.Cells.Item(17, 9).Value = 1.043644113806419
$ws.Cells.Item(17, 10).Value = 1.047530813257024
$ws.Cells.Item(17, 11).Value = 1.052215575636006
$ws.Cells.Item(17, 12).Value = 1.058199905924293
$ws.Cells.Item(17, 13).Value = 1.065075111299919
$ws.Cells.Item(17, 14).Value = 1.019721556162823

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.041566951316141
$ws.Cells.Item(18, 4).Value = 1.049033534591311
$ws.Cells.Item(18, 5).Value = 1.055068667139839
$ws.Cells.Item(18, 6).Value = 1.061967159268086
$ws.Cells.Item(18, 9).Value = 1.043685907891172
$ws.Cells.Item(18, 10).Value = 1.047633015381746
$ws.Cells.Item(18, 11).Value = 1.052313259047174
$ws.Cells.Item(18, 12).Value = 1.058328326391616
$ws.Cells.Item(18, 13).Value = 1.065204189079848
$ws.Cells.Item(18, 14).Value = 1.01975586596208

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.041628954845992
$ws.Cells.Item(19, 4).Value = 1.049081105687486
$ws.Cells.Item(19, 5).Value = 1.055126668561639
$ws.Cells.Item(19, 6).Value = 1.062025278766703
$ws.Cells.Item(19, 9).Value = 1.043700144005216
$ws.Cells.Item(19, 10).Value = 1.047667866742897
$ws.Cells.Item(19, 11).Value = 1.052346567631942
$ws.Cells.Item(19, 12).Value = 1.05837212444775
$ws.Cells.Item(19, 13).Value = 1.065248210588807
$ws.Cells.Item(19, 14).Value = 1.019767565045706

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.041351728538235
$ws.Cells.Item(20, 4).Value = 1.048868416662814
$ws.Cells.Item(20, 5).Value = 1.054867365563511
$ws.Cells.Item(20, 6).Value = 1.061765445076709
$ws.Cells.Item(20, 9).Value = 1.043636419180118
$ws.Cells.Item(20, 10).Value = 1.047512015415961
$ws.Cells.Item(20, 11).Value = 1.052197608023591
$ws.Cells.Item(20, 12).Value = 1.058176288702049
$ws.Cells.Item(20, 13).Value = 1.065051372850482
$ws.Cells.Item(20, 14).Value = 1.019715245293083

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.040451642266896
$ws.Cells.Item(21, 4).Value = 1.048178007840733
$ws.Cells.Item(21, 5).Value = 1.05402598851845
$ws.Cells.Item(21, 6).Value = 1.060922296805701
$ws.Cells.Item(21, 9).Value = 1.043428236182922
$ws.Cells.Item(21, 10).Value = 1.047005547164105
$ws.Cells.Item(21, 11).Value = 1.051713406116378
$ws.Cells.Item(21, 12).Value = 1.057540307980292
$ws.Cells.Item(21, 13).Value = 1.064412088062554
$ws.Cells.Item(21, 14).Value = 1.019545173770589

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039886654774484
$ws.Cells.Item(22, 4).Value = 1.047744746140969
$ws.Cells.Item(22, 5).Value = 1.053498251333761
$ws.Cells.Item(22, 6).Value = 1.060393411023491
$ws.Cells.Item(22, 9).Value = 1.043296575173306
$ws.Cells.Item(22, 10).Value = 1.046687286054953
$ws.Cells.Item(22, 11).Value = 1.051409038093473
$ws.Cells.Item(22, 12).Value = 1.057140987381434
$ws.Cells.Item(22, 13).Value = 1.064010654663164
$ws.Cells.Item(22, 14).Value = 1.019438264944181

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.04018607628836
$ws.Cells.Item(23, 4).Value = 1.047974347426754
$ws.Cells.Item(23, 5).Value = 1.053777893692339
$ws.Cells.Item(23, 6).Value = 1.060673665574583
$ws.Cells.Item(23, 9).Value = 1.043366443706119
$ws.Cells.Item(23, 10).Value = 1.046855985171739
$ws.Cells.Item(23, 11).Value = 1.051570382323346
$ws.Cells.Item(23, 12).Value = 1.05735262247581
$ws.Cells.Item(23, 13).Value = 1.064223413133809
$ws.Cells.Item(23, 14).Value = 1.019494937087046

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.041366834265733
$ws.Cells.Item(24, 4).Value = 1.04888000530775
$ws.Cells.Item(24, 5).Value = 1.05488149272779
$ws.Cells.Item(24, 6).Value = 1.061779601340036
$ws.Cells.Item(24, 9).Value = 1.043639896315955
$ws.Cells.Item(24, 10).Value = 1.047520509294529
$ws.Cells.Item(24, 11).Value = 1.052205726794208
$ws.Cells.Item(24, 12).Value = 1.058186960125375
$ws.Cells.Item(24, 13).Value = 1.065062099063414
$ws.Cells.Item(24, 14).Value = 1.019718096897263

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.042741032299345
$ws.Cells.Item(25, 4).Value = 1.049934498950134
$ws.Cells.Item(25, 5).Value = 1.056167614219399
$ws.Cells.Item(25, 6).Value = 1.063068280815149
$ws.Cells.Item(25, 9).Value = 1.043953851569527
$ws.Cells.Item(25, 10).Value = 1.048292369230241
$ws.Cells.Item(25, 11).Value = 1.052943260166687
$ws.Cells.Item(25, 12).Value = 1.059157480182642
$ws.Cells.Item(25, 13).Value = 1.066037509383267
$ws.Cells.Item(25, 14).Value = 1.019977139242566

Write-Output "Updated vm_pu values for rows 2-25 (380 kV case)"